$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.426.37"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "2.015.25"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'259.42"
$ws.Range("E5").Value = "  +5.42%  "

$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'56.12"
$ws.Range("E8").Value = "  -6.11%  "

$ws.Range("E9").Value = "  +1.20%  "

$ws.Range("E10").Value = "  -3.84%  "

$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("D12").Value = "'14.31"
$ws.Range("E12").Value = "  -4.86%  "

$ws.Range("D13").Value = "2.311.88"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").Value = "'0.806"
$ws.Range("E14").Value = "  -4.58%  "

$ws.Range("D15").Value = "'20.82"
$ws.Range("E15").Value = "  -7.25%  "

$ws.Range("D16").Value = "'5.27"
$ws.Range("E16").Value = "  -2.91%  "

$ws.Range("D17").Value = "2.013.96"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").Value = "37.301.48"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "'69.76"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "0.0₃0843"
$ws.Range("E20").Value = "  -2.38%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "'228.39"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").Value = "'2.67"
$ws.Range("E23").Value = "  +8.27%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  -1.41%  "

$ws.Range("D26").Value = "'164.47"
$ws.Range("E26").Value = "  +0.59%  "

$ws.Range("D27").Value = "'9.00"
$ws.Range("E27").Value = "  -4.57%  "

$ws.Range("D28").Value = "'19.71"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("E29").Value = "  -9.19%  "

$ws.Range("D30").Value = "'1.32"
$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("E32").Value = "  -3.09%  "

$ws.Range("D33").Value = "'0.0650"
$ws.Range("E33").Value = "  -0.83%  "

$ws.Range("D34").Value = "'4.54"
$ws.Range("E34").Value = "  +1.16%  "

$ws.Range("D35").Value = "'2.40"
$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("D36").Value = "'1.82"
$ws.Range("E36").Value = "  +1.02%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").Value = "'3.35"
$ws.Range("E38").Value = "  +1.51%  "

$ws.Range("D39").Value = "'5.19"
$ws.Range("E39").Value = "  -3.57%  "

$ws.Range("D40").Value = "'3.03"
$ws.Range("E40").Value = "  +3.65%  "

$ws.Range("E41").Value = "  +3.32%  "

$ws.Range("D42").Value = "'0.0941"
$ws.Range("E42").Value = "  -3.88%  "

$ws.Range("D43").Value = "'0.0213"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("D44").Value = "1.393.30"
$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("D45").Value = "'90.36"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("D46").Value = "'15.69"
$ws.Range("E46").Value = "  -5.77%  "

$ws.Range("D47").Value = "'1.02"
$ws.Range("E47").Value = "  -2.06%  "

$ws.Range("D48").Value = "'7.10"
$ws.Range("E48").Value = "  -4.61%  "

$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("D50").Value = "2.203.40"
$ws.Range("E50").Value = "  +0.65%  "

$ws.Range("D51").Value = "'1.95"
$ws.Range("E51").Value = "  -4.54%  "
